# Update "想去人数" (want-to-go count) figures in column F
# on the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F2"  = 252
    "F3"  = 276
    "F7"  = 6858
    "F11" = 90
    "F16" = 238
    "F17" = 597
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
